# Update column G ("K") values for rows 2-11 on the active worksheet.
# This mirrors the commit "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" which only changed the numeric
# values stored in column G for data rows 2 through 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 8
    3  = 4
    4  = 1
    5  = 3
    6  = 2
    7  = 4
    8  = 0
    9  = 4
    10 = 1
    11 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
